$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" column header in F1, matching the style/formatting
# of the existing header cells (e.g. E1: bold font, borders, centered).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the time_taken values for each data row as plain text strings.
$ws.Range("F2").Value = "2021-10-05 13:40:16.693406"
$ws.Range("F3").Value = "2021-10-05 13:40:16.693418"
$ws.Range("F4").Value = "2021-10-05 13:40:16.693422"
